$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 130.955829
$ws.Range("H2").Value = 392.867487
$ws.Range("I2").Value = 0.5336535908353144
$ws.Range("J2").Value = 0.5336535908353144
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 28.72417333333333
$ws.Range("N2").Value = 86.17251999999999
$ws.Range("O2").Value = 0.4233259107972328
$ws.Range("P2").Value = 0.4233259107972328
$ws.Range("Q2").Value = 3761.597931206359
$ws.Range("R2").Value = 33854.38138085724
$ws.Range("S2").Value = 0.2259093923905733
$ws.Range("T2").Value = 0.2259093923905733

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 130.955829
$ws.Range("H3").Value = 392.867487
$ws.Range("I3").Value = 0.5336535908353144
$ws.Range("J3").Value = 0.5336535908353144
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 30.56986233333333
$ws.Range("N3").Value = 91.709587
$ws.Range("O3").Value = 0.4505269713084062
$ws.Range("P3").Value = 0.4505269713084062
$ws.Range("Q3").Value = 4003.301664277541
$ws.Range("R3").Value = 36029.71497849787
$ws.Range("S3").Value = 0.2404253360068896
$ws.Range("T3").Value = 0.2404253360068896

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 130.955829
$ws.Range("H4").Value = 392.867487
$ws.Range("I4").Value = 0.5336535908353144
$ws.Range("J4").Value = 0.5336535908353144
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.559531999999999
$ws.Range("N4").Value = 25.678596
$ws.Range("O4").Value = 0.126147117894361
$ws.Range("P4").Value = 0.126147117894361
$ws.Range("Q4").Value = 1120.920608912028
$ws.Range("R4").Value = 10088.28548020825
$ws.Range("S4").Value = 0.06731886243785148
$ws.Range("T4").Value = 0.06731886243785148

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 66.39541
$ws.Range("H5").Value = 199.18623
$ws.Range("I5").Value = 0.2705656497465488
$ws.Range("J5").Value = 0.2705656497465488
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 28.72417333333333
$ws.Range("N5").Value = 86.17251999999999
$ws.Range("O5").Value = 0.4233259107972328
$ws.Range("P5").Value = 0.4233259107972328
$ws.Range("Q5").Value = 1907.153265377733
$ws.Range("R5").Value = 17164.3793883996
$ws.Range("S5").Value = 0.1145374501094029
$ws.Range("T5").Value = 0.1145374501094029

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 66.39541
$ws.Range("H6").Value = 199.18623
$ws.Range("I6").Value = 0.2705656497465488
$ws.Range("J6").Value = 0.2705656497465488
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.56986233333333
$ws.Range("N6").Value = 91.709587
$ws.Range("O6").Value = 0.4505269713084062
$ws.Range("P6").Value = 0.4505269713084062
$ws.Range("Q6").Value = 2029.698543265223
$ws.Range("R6").Value = 18267.28688938701
$ws.Range("S6").Value = 0.1218971227204037
$ws.Range("T6").Value = 0.1218971227204037

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 66.39541
$ws.Range("H7").Value = 199.18623
$ws.Range("I7").Value = 0.2705656497465488
$ws.Range("J7").Value = 0.2705656497465488
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.559531999999999
$ws.Range("N7").Value = 25.678596
$ws.Range("O7").Value = 0.126147117894361
$ws.Range("P7").Value = 0.126147117894361
$ws.Range("Q7").Value = 568.3136365481199
$ws.Range("R7").Value = 5114.82272893308
$ws.Range("S7").Value = 0.03413107691674228
$ws.Range("T7").Value = 0.03413107691674228

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 48.043585
$ws.Range("H8").Value = 144.130755
$ws.Range("I8").Value = 0.1957807594181367
$ws.Range("J8").Value = 0.1957807594181367
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 28.72417333333333
$ws.Range("N8").Value = 86.17251999999999
$ws.Range("O8").Value = 0.4233259107972328
$ws.Range("P8").Value = 0.4233259107972328
$ws.Range("Q8").Value = 1380.012263094733
$ws.Range("R8").Value = 12420.1103678526
$ws.Range("S8").Value = 0.08287906829725664
$ws.Range("T8").Value = 0.08287906829725664

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 48.043585
$ws.Range("H9").Value = 144.130755
$ws.Range("I9").Value = 0.1957807594181367
$ws.Range("J9").Value = 0.1957807594181367
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.56986233333333
$ws.Range("N9").Value = 91.709587
$ws.Range("O9").Value = 0.4505269713084062
$ws.Range("P9").Value = 0.4505269713084062
$ws.Range("Q9").Value = 1468.685779449798
$ws.Range("R9").Value = 13218.17201504819
$ws.Range("S9").Value = 0.08820451258111286
$ws.Range("T9").Value = 0.08820451258111285

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 48.043585
$ws.Range("H10").Value = 144.130755
$ws.Range("I10").Value = 0.1957807594181367
$ws.Range("J10").Value = 0.1957807594181367
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.559531999999999
$ws.Range("N10").Value = 25.678596
$ws.Range("O10").Value = 0.126147117894361
$ws.Range("P10").Value = 0.126147117894361
$ws.Range("Q10").Value = 411.23060320222
$ws.Range("R10").Value = 3701.07542881998
$ws.Range("S10").Value = 0.02469717853976722
$ws.Range("T10").Value = 0.02469717853976721
